$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 1136
$wsExhibition.Range("F6").Value = 9
$wsExhibition.Range("F8").Value = 255
$wsExhibition.Range("F14").Value = 161
$wsExhibition.Range("F15").Value = 12928
$wsExhibition.Range("F17").Value = 5308
$wsExhibition.Range("F18").Value = 5538

# --- Sheet "演出" (Performance) ---
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 146

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1136
$wsAll.Range("F6").Value = 9
$wsAll.Range("F8").Value = 255
$wsAll.Range("F14").Value = 161
$wsAll.Range("F15").Value = 12928
$wsAll.Range("F16").Value = 146
$wsAll.Range("F19").Value = 5308
$wsAll.Range("F20").Value = 5538
